$p = $ppt.ActivePresentation
$p.Slides.Add(19, 1)
